$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "StatQuery" column is being added between the existing "query" (A)
# and "dbExcel" (old B, now C) columns - inserting a whole column shifts
# the old B->C and old C->D, which matches the target layout.
$ws.Columns("B").Insert()

# The new column B should carry the same (wide) width as column A.
# Excel's ColumnWidth setter only offers integer-character precision in
# this engine, so 75 is the closest achievable value to column A's
# 75.81640625 stored width.
$ws.Columns("B").ColumnWidth = 75

# New header cell for the inserted column.
$ws.Range("B1").Value = "StatQuery"

# New long Cypher "stat" query that accompanies the header - goes in row 2
# of the new column, matching the wrap style already used by A2.
$statQuery = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.race IN  ['NOT_REPORTED']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").Value = $statQuery

# Note: B2 already inherits the wrap-text formatting (matching A2/row 2's
# long-query cells) from the Columns("B").Insert() above - it copies the
# left-hand column's cell formatting for the new column automatically.
